$p = $ppt.ActivePresentation

# --- Slide 1 ("Title Slide"): fill in the title & subtitle text ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Supervisory Meeting"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "09/02/24"

# --- Add a new slide 2 using the "Title and Content" layout ---
$s2 = $p.Slides.Add(2, 2)

# Title placeholder
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Agenda`t"

# Body / content placeholder: one bullet per line, with a final
# non-bulleted blank line.
$body = $s2.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Elevator Pitch"
[void]$body.InsertAfter("`rKeywords Search")
[void]$body.InsertAfter("`rSearch Results`t")
[void]$body.InsertAfter("`r")
$body.Paragraphs(4).ParagraphFormat.Bullet.Type = 0
